$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1838.8334
$ws.Range("I39").Value = 1345.875
$ws.Range("K39").Value = 4037.625
$ws.Range("M39").Value = -3741.625
$ws.Range("H40").Value = 3425
$ws.Range("I40").Value = 3425
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3425
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3250
$ws.Range("N40").ClearContents()
$ws.Range("H53").Value = 3213
$ws.Range("I53").Value = 3027.6365
$ws.Range("J53").Value = 3416.9
$ws.Range("K53").Value = 3027.6365
$ws.Range("L53").Value = 3416.9
$ws.Range("M53").Value = -2390.6365
$ws.Range("N53").Value = -4690.9
$ws.Range("H62").Value = 62512460
$ws.Range("I62").Value = 142858300
$ws.Range("K62").Value = 142858300
$ws.Range("M62").Value = -142857676
$ws.Range("H64").Value = 7831.6665
$ws.Range("I64").Value = 7750
$ws.Range("K64").Value = 7750
$ws.Range("M64").Value = -7502
$ws.Range("H65").Value = 62512460
$ws.Range("I65").Value = 142858300
$ws.Range("K65").Value = 714291500
$ws.Range("M65").Value = -714288380
$ws.Range("H67").Value = 7831.6665
$ws.Range("I67").Value = 7750
$ws.Range("K67").Value = 7750
$ws.Range("M67").Value = -6892
$ws.Range("H132").Value = 2209.611
$ws.Range("I132").Value = 2209.611
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6628.833
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4098.833
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 2043.9395
$ws.Range("I137").Value = 1998.3871
$ws.Range("K137").Value = 5995.1613
$ws.Range("M137").Value = -3445.1613
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2930.625
$ws.Range("I63").Value = 2930.625
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2930.625
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2244.625
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 2930.625
$ws.Range("I66").Value = 2930.625
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 14653.125
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -11221.125
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 17834.71
$ws.Range("I74").Value = 24365.54
$ws.Range("K74").Value = 24365.54
$ws.Range("M74").Value = -23491.54
$ws.Range("H77").Value = 17834.71
$ws.Range("I77").Value = 24365.54
$ws.Range("K77").Value = 121827.7
$ws.Range("M77").Value = -117459.7
$ws.Range("H132").Value = 5764.1777
$ws.Range("I132").Value = 4222.0713
$ws.Range("J132").Value = 8304.117
$ws.Range("K132").Value = 12666.2139
$ws.Range("L132").Value = 24912.351
$ws.Range("M132").Value = -10136.2139
$ws.Range("N132").Value = -29972.351
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2269.4849
$ws.Range("I94").Value = 577.45
$ws.Range("J94").Value = 4872.615
$ws.Range("K94").Value = 577.45
$ws.Range("L94").Value = 4872.615
$ws.Range("M94").Value = -126.45
$ws.Range("N94").Value = -5774.615
$ws.Range("H99").Value = 3639336.8
$ws.Range("J99").Value = 10104035
$ws.Range("L99").Value = 10104035
$ws.Range("N99").Value = -10107031
$ws.Range("H128").Value = 3318.111
$ws.Range("I128").Value = 3318.111
$ws.Range("K128").Value = 9954.332999999999
$ws.Range("M128").Value = -7464.332999999999
$ws.Range("H134").Value = 5824.273
$ws.Range("I134").Value = 2078.353
$ws.Range("J134").Value = 8182.815
$ws.Range("K134").Value = 6235.059
$ws.Range("L134").Value = 24548.445
$ws.Range("M134").Value = -3700.059
$ws.Range("N134").Value = -29618.445
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 239.66667
$ws.Range("I22").Value = 220
$ws.Range("J22").Value = 279
$ws.Range("K22").Value = 220
$ws.Range("L22").Value = 279
$ws.Range("M22").Value = 130
$ws.Range("N22").Value = -979
$ws.Range("H31").Value = 6419.66
$ws.Range("I31").Value = 2403.0605
$ws.Range("K31").Value = 2403.0605
$ws.Range("M31").Value = -2108.0605
$ws.Range("H34").Value = 6419.66
$ws.Range("I34").Value = 2403.0605
$ws.Range("K34").Value = 2403.0605
$ws.Range("M34").Value = -2201.0605
$ws.Range("H58").Value = 8933087
$ws.Range("I58").Value = 13159227
$ws.Range("K58").Value = 13159227
$ws.Range("M58").Value = -13159024
$ws.Range("H94").Value = 645.15
$ws.Range("I94").Value = 788.3333
$ws.Range("J94").Value = 583.7857
$ws.Range("K94").Value = 788.3333
$ws.Range("L94").Value = 583.7857
$ws.Range("M94").Value = -337.3333
$ws.Range("N94").Value = -1485.7857
$ws.Range("H134").Value = 5855.8423
$ws.Range("I134").Value = 1637.8948
$ws.Range("J134").Value = 10073.789
$ws.Range("K134").Value = 4913.6844
$ws.Range("L134").Value = 30221.367
$ws.Range("M134").Value = -2378.6844
$ws.Range("N134").Value = -35291.367
$ws.Range("H136").Value = 8933087
$ws.Range("I136").Value = 13159227
$ws.Range("K136").Value = 39477681
$ws.Range("M136").Value = -39475131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 12821405
$ws.Range("I14").Value = 12821405
$ws.Range("K14").Value = 38464215
$ws.Range("M14").Value = -38464042
$ws.Range("H33").Value = 23809704
$ws.Range("I33").Value = 55555596
$ws.Range("J33").Value = 286.75
$ws.Range("K33").Value = 333333576
$ws.Range("L33").Value = 1720.5
$ws.Range("M33").Value = -333333293
$ws.Range("N33").Value = -2286.5
$ws.Range("H46").Value = 101220.2
$ws.Range("I46").Value = 333567.34
$ws.Range("J46").Value = 1642.8572
$ws.Range("K46").Value = 1000702.02
$ws.Range("L46").Value = 4928.571599999999
$ws.Range("M46").Value = -1000611.02
$ws.Range("N46").Value = -5110.571599999999
$ws.Range("H69").Value = 4000
$ws.Range("I69").Value = 2500
$ws.Range("K69").Value = 7500
$ws.Range("M69").Value = -6689
$ws.Range("H72").Value = 4000
$ws.Range("I72").Value = 2500
$ws.Range("K72").Value = 22500
$ws.Range("M72").Value = -18444
$ws.Range("H140").Value = 2032.6666
$ws.Range("I140").Value = 1374.1666
$ws.Range("K140").Value = 4122.4998
$ws.Range("M140").Value = 1057.5002
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2277.3462
$ws.Range("I80").Value = 2156.4167
$ws.Range("J80").Value = 2381
$ws.Range("K80").Value = 2156.4167
$ws.Range("L80").Value = 2381
$ws.Range("M80").Value = -1158.4167
$ws.Range("N80").Value = -4377
$ws.Range("H83").Value = 2277.3462
$ws.Range("I83").Value = 2156.4167
$ws.Range("J83").Value = 2381
$ws.Range("K83").Value = 10782.0835
$ws.Range("L83").Value = 11905
$ws.Range("M83").Value = -5790.083500000001
$ws.Range("N83").Value = -21889
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3817.2068
$ws.Range("I61").Value = 1097.5454
$ws.Range("K61").Value = 1097.5454
$ws.Range("M61").Value = -895.5454
$ws.Range("H113").Value = 3817.2068
$ws.Range("I113").Value = 1097.5454
$ws.Range("K113").Value = 1097.5454
$ws.Range("M113").Value = 1072.4546
$ws.Range("H132").Value = 7697177
$ws.Range("I132").Value = 13890594
$ws.Range("K132").Value = 41671782
$ws.Range("M132").Value = -41669252
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 76013.28999999999
$ws.Range("J62").Value = 4599
$ws.Range("L62").Value = 4599
$ws.Range("N62").Value = -5847
$ws.Range("H65").Value = 76013.28999999999
$ws.Range("J65").Value = 4599
$ws.Range("L65").Value = 22995
$ws.Range("N65").Value = -29235
$ws.Range("H107").Value = 1259.4375
$ws.Range("I107").Value = 1230.6666
$ws.Range("J107").Value = 1296.4286
$ws.Range("K107").Value = 3691.9998
$ws.Range("L107").Value = 3889.2858
$ws.Range("M107").Value = -1771.9998
$ws.Range("N107").Value = -7729.2858
$ws.Range("H122").Value = 233517.06
$ws.Range("I122").Value = 822681.6
$ws.Range("J122").Value = 6915.3076
$ws.Range("K122").Value = 2468044.8
$ws.Range("L122").Value = 20745.9228
$ws.Range("M122").Value = -2465594.8
$ws.Range("N122").Value = -25645.9228
$ws.Range("H132").Value = 8204035.5
$ws.Range("I132").Value = 10002910
$ws.Range("J132").Value = 27330.455
$ws.Range("K132").Value = 30008730
$ws.Range("L132").Value = 81991.36500000001
$ws.Range("M132").Value = -30006200
$ws.Range("N132").Value = -87051.36500000001
